$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Cadm3"
$ws.Cells.Item(2,3).Value = "Cadm3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 3.443291333333333
$ws.Cells.Item(2,8).Value = 10.329874
$ws.Cells.Item(2,9).Value = 0.07497468122035157
$ws.Cells.Item(2,10).Value = 0.07497468122035157
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 3.443291333333333
$ws.Cells.Item(2,14).Value = 10.329874
$ws.Cells.Item(2,15).Value = 0.07497468122035157
$ws.Cells.Item(2,16).Value = 0.07497468122035157
$ws.Cells.Item(2,17).Value = 11.85625520620844
$ws.Cells.Item(2,18).Value = 106.706296855876
$ws.Cells.Item(2,19).Value = 0.005621202824093338
$ws.Cells.Item(2,20).Value = 0.005621202824093338

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Cadm3"
$ws.Cells.Item(3,3).Value = "Cadm3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 3.443291333333333
$ws.Cells.Item(3,8).Value = 10.329874
$ws.Cells.Item(3,9).Value = 0.07497468122035157
$ws.Cells.Item(3,10).Value = 0.07497468122035157
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 39.68460733333333
$ws.Cells.Item(3,14).Value = 119.053822
$ws.Cells.Item(3,15).Value = 0.8640978924345524
$ws.Cells.Item(3,16).Value = 0.8640978924345523
$ws.Cells.Item(3,17).Value = 136.6456644976031
$ws.Cells.Item(3,18).Value = 1229.810980478428
$ws.Cells.Item(3,19).Value = 0.06478546402845821
$ws.Cells.Item(3,20).Value = 0.0647854640284582

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Cadm3"
$ws.Cells.Item(4,3).Value = "Cadm3"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 3.443291333333333
$ws.Cells.Item(4,8).Value = 10.329874
$ws.Cells.Item(4,9).Value = 0.07497468122035157
$ws.Cells.Item(4,10).Value = 0.07497468122035157
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.798156333333333
$ws.Cells.Item(4,14).Value = 8.394469
$ws.Cells.Item(4,15).Value = 0.06092742634509613
$ws.Cells.Item(4,16).Value = 0.06092742634509613
$ws.Cells.Item(4,17).Value = 9.634867451878444
$ws.Cells.Item(4,18).Value = 86.71380706690601
$ws.Cells.Item(4,19).Value = 0.004568014367800032
$ws.Cells.Item(4,20).Value = 0.004568014367800032

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Cadm3"
$ws.Cells.Item(5,3).Value = "Cadm3"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 39.68460733333333
$ws.Cells.Item(5,8).Value = 119.053822
$ws.Cells.Item(5,9).Value = 0.8640978924345524
$ws.Cells.Item(5,10).Value = 0.8640978924345523
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 3.443291333333333
$ws.Cells.Item(5,14).Value = 10.329874
$ws.Cells.Item(5,15).Value = 0.07497468122035157
$ws.Cells.Item(5,16).Value = 0.07497468122035157
$ws.Cells.Item(5,17).Value = 136.6456644976031
$ws.Cells.Item(5,18).Value = 1229.810980478428
$ws.Cells.Item(5,19).Value = 0.06478546402845821
$ws.Cells.Item(5,20).Value = 0.0647854640284582

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Cadm3"
$ws.Cells.Item(6,3).Value = "Cadm3"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 39.68460733333333
$ws.Cells.Item(6,8).Value = 119.053822
$ws.Cells.Item(6,9).Value = 0.8640978924345524
$ws.Cells.Item(6,10).Value = 0.8640978924345523
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 39.68460733333333
$ws.Cells.Item(6,14).Value = 119.053822
$ws.Cells.Item(6,15).Value = 0.8640978924345524
$ws.Cells.Item(6,16).Value = 0.8640978924345523
$ws.Cells.Item(6,17).Value = 1574.868059200854
$ws.Cells.Item(6,18).Value = 14173.81253280768
$ws.Cells.Item(6,19).Value = 0.7466651677098353
$ws.Cells.Item(6,20).Value = 0.746665167709835

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Cadm3"
$ws.Cells.Item(7,3).Value = "Cadm3"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 39.68460733333333
$ws.Cells.Item(7,8).Value = 119.053822
$ws.Cells.Item(7,9).Value = 0.8640978924345524
$ws.Cells.Item(7,10).Value = 0.8640978924345523
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.798156333333333
$ws.Cells.Item(7,14).Value = 8.394469
$ws.Cells.Item(7,15).Value = 0.06092742634509613
$ws.Cells.Item(7,16).Value = 0.06092742634509613
$ws.Cells.Item(7,17).Value = 111.0437353456131
$ws.Cells.Item(7,18).Value = 999.3936181105181
$ws.Cells.Item(7,19).Value = 0.05264726069625898
$ws.Cells.Item(7,20).Value = 0.05264726069625898

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Cadm3"
$ws.Cells.Item(8,3).Value = "Cadm3"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.798156333333333
$ws.Cells.Item(8,8).Value = 8.394469
$ws.Cells.Item(8,9).Value = 0.06092742634509613
$ws.Cells.Item(8,10).Value = 0.06092742634509613
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 3.443291333333333
$ws.Cells.Item(8,14).Value = 10.329874
$ws.Cells.Item(8,15).Value = 0.07497468122035157
$ws.Cells.Item(8,16).Value = 0.07497468122035157
$ws.Cells.Item(8,17).Value = 9.634867451878444
$ws.Cells.Item(8,18).Value = 86.71380706690601
$ws.Cells.Item(8,19).Value = 0.004568014367800032
$ws.Cells.Item(8,20).Value = 0.004568014367800032

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Cadm3"
$ws.Cells.Item(9,3).Value = "Cadm3"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.798156333333333
$ws.Cells.Item(9,8).Value = 8.394469
$ws.Cells.Item(9,9).Value = 0.06092742634509613
$ws.Cells.Item(9,10).Value = 0.06092742634509613
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 39.68460733333333
$ws.Cells.Item(9,14).Value = 119.053822
$ws.Cells.Item(9,15).Value = 0.8640978924345524
$ws.Cells.Item(9,16).Value = 0.8640978924345523
$ws.Cells.Item(9,17).Value = 111.0437353456131
$ws.Cells.Item(9,18).Value = 999.3936181105181
$ws.Cells.Item(9,19).Value = 0.05264726069625898
$ws.Cells.Item(9,20).Value = 0.05264726069625898

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Cadm3"
$ws.Cells.Item(10,3).Value = "Cadm3"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 2.798156333333333
$ws.Cells.Item(10,8).Value = 8.394469
$ws.Cells.Item(10,9).Value = 0.06092742634509613
$ws.Cells.Item(10,10).Value = 0.06092742634509613
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.798156333333333
$ws.Cells.Item(10,14).Value = 8.394469
$ws.Cells.Item(10,15).Value = 0.06092742634509613
$ws.Cells.Item(10,16).Value = 0.06092742634509613
$ws.Cells.Item(10,17).Value = 7.829678865773445
$ws.Cells.Item(10,18).Value = 70.46710979196101
$ws.Cells.Item(10,19).Value = 0.003712151281037113
$ws.Cells.Item(10,20).Value = 0.003712151281037113

